# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# The underlying data rows got re-ordered / corrected. Column A (the
# sequential "id" row index) stays put; every other column (B..AC) for
# the affected rows is exchanged between rows as described below:
#   16 <-> 17
#   63 <-> 64
#   85 <-> 86
#   98 <-> 99
#   101 -> gets old 104, 103 -> gets old 101, 104 -> gets old 103 (3-way rotation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

function Get-RowValues([int]$row) {
    $vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues([int]$row, $vals) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$c]
    }
}

function Swap-Rows([int]$rowA, [int]$rowB) {
    $valsA = Get-RowValues $rowA
    $valsB = Get-RowValues $rowB
    Set-RowValues $rowA $valsB
    Set-RowValues $rowB $valsA
}

# Simple pairwise swaps
Swap-Rows 16 17
Swap-Rows 63 64
Swap-Rows 85 86
Swap-Rows 98 99

# 3-way rotation among rows 101, 103, 104:
#   new(101) = old(104)
#   new(103) = old(101)
#   new(104) = old(103)
$vals101 = Get-RowValues 101
$vals103 = Get-RowValues 103
$vals104 = Get-RowValues 104

Set-RowValues 101 $vals104
Set-RowValues 103 $vals101
Set-RowValues 104 $vals103
